# "fully operating charging station" - add a battery (bat1) node to the
# power-flow matrix on Sheet1: a new "From" column (D) and a new "To" row
# (row 4, between net1 and charging_station1), wired into every other node.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column headers (row 1): From net1 / pv1 / bat1 -------------------------
$ws.Range("B1").Value = "P_from_net1"
$ws.Range("C1").Value = "P_from_pv1"
$ws.Range("D1").Value = "P_from_bat1"

# --- Row labels (column A): To demand1 / net1 / bat1 / charging_station1 ----
$ws.Range("A2").Value = "P_to_demand1"
$ws.Range("A3").Value = "P_to_net1"
$ws.Range("A4").Value = "P_to_bat1"
$ws.Range("A5").Value = "P_to_charging_station1"

# --- Row 2: -> demand1 --------------------------------------------------
$ws.Range("B2").Value = "P_net1_demand1"
$ws.Range("C2").Value = "P_pv1_demand1"
$ws.Range("D2").Value = "P_bat1_demand1"

# --- Row 3: -> net1 (net1->net1 stays the numeric placeholder 0) -------
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "P_pv1_net1"
$ws.Range("D3").Value = "P_bat1_net1"

# --- Row 4: -> bat1 (new row; bat1->bat1 is the numeric placeholder 0) -
$ws.Range("B4").Value = "P_net1_bat1"
$ws.Range("C4").Value = "P_pv1_bat1"
$ws.Range("D4").Value = 0

# --- Row 5: -> charging_station1 (was row 4) ----------------------------
$ws.Range("B5").Value = "P_net1_charging_station1"
$ws.Range("C5").Value = "P_pv1_charging_station1"
$ws.Range("D5").Value = "P_bat1_charging_station1"

# --- Carry the header style (bold, centered, bordered - style index 1 in
#     the original file) onto the brand-new D1/A5 header cells; every other
#     touched cell already had that style (or is a plain data cell) and
#     keeps it automatically. -------------------------------------------
$ws.Range("B1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null

$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
